$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExternalCreds")

# Update header cell A1 text
$ws.Range("A1").Value = "RowSelection"

# Reset the sheet's selection back to A1 (removes the stored F13 selection)
$ws.Activate()
$ws.Range("A1").Select()
